$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 71, shifting existing rows 71:196 down to 72:197.
$ws.Rows(71).Insert()

# Populate the newly inserted row 71 with the new record's data.
$ws.Range("A71").Value = 7
$ws.Range("B71").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C71").Value = "Ñuble"
$ws.Range("D71").Value = 44533
$ws.Range("E71").Value = 16
$ws.Range("F71").Value = 100112008
$ws.Range("G71").Value = "Coliflor"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 300
$ws.Range("K71").Value = 700
$ws.Range("L71").Value = 800
$ws.Range("M71").Value = 750
$ws.Range("N71").Value = "$/unidad"
$ws.Range("O71").Value = "Región del Maule"
$ws.Range("P71").Value = 750
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = "Hortaliza"
